# Convention change to support multi-axle vehicles.
# Renames the "front axle" instance name / class label used in every
# Body_1Axle sheet of this workbook:
#   Instance "sAxleF"    -> "sAxle1"
#   class    "Body_1Axle" -> "Body_Axle1"
# Also nudges the remembered cell selection on each sheet (cosmetic,
# matches what Excel records after the edit/save).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Trailer_Elula", "Trailer_Elula_Unstable", "Trailer_Thwala")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # A5 = "Instance" row value, H4 = "class" row value
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("H4").Value = "Body_Axle1"
}

# Restore the per-sheet "last selected cell" the way the author's Excel
# session left them.
$wsElula = $wb.Worksheets.Item("Trailer_Elula")
$wsElula.Range("M7").Select()

$wsUnstable = $wb.Worksheets.Item("Trailer_Elula_Unstable")
$wsUnstable.Range("H4").Select()

$wsThwala = $wb.Worksheets.Item("Trailer_Thwala")
$wsThwala.Range("H4").Select()

# Keep the originally active sheet/tab selected.
$wsElula.Activate()
